$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 5: SK / SK증권제12호스팩 ---
$ws.Rows.Item(5).Insert()

$ws.Range("B5").NumberFormat = "@"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("G5").NumberFormat = "@"

$ws.Cells.Item(5,1).Value = "SK"
$ws.Cells.Item(5,2).Value = "2024-04-23"
$ws.Cells.Item(5,3).Value = "SK증권제12호스팩"
$ws.Cells.Item(5,4).Value = "SK"
$ws.Cells.Item(5,5).Value = "SK"
$ws.Cells.Item(5,6).Value = "2024-04-26"
$ws.Cells.Item(5,7).Value = "2024-05-07"
$ws.Cells.Item(5,8).Value = 6000
$ws.Cells.Item(5,9).Value = 3000000
$ws.Cells.Item(5,10).Value = 2000
$ws.Cells.Item(5,11).Value = 0
$ws.Cells.Item(5,12).Value = 100

$ws.Range("B5").ClearFormats()
$ws.Range("F5").ClearFormats()
$ws.Range("G5").ClearFormats()

# --- Insert new row 12: 한국 / 코칩 ---
$ws.Rows.Item(12).Insert()

$ws.Range("B12").NumberFormat = "@"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("G12").NumberFormat = "@"

$ws.Cells.Item(12,1).Value = "한국"
$ws.Cells.Item(12,2).Value = "2024-04-24"
$ws.Cells.Item(12,3).Value = "코칩"
$ws.Cells.Item(12,4).Value = "한국"
$ws.Cells.Item(12,5).Value = "한국"
$ws.Cells.Item(12,6).Value = "2024-04-29"
$ws.Cells.Item(12,7).Value = "2024-05-07"
$ws.Cells.Item(12,8).Value = 27000
$ws.Cells.Item(12,9).Value = 1500000
$ws.Cells.Item(12,10).Value = 18000
$ws.Cells.Item(12,11).Value = 0
$ws.Cells.Item(12,12).Value = 100

$ws.Range("B12").ClearFormats()
$ws.Range("F12").ClearFormats()
$ws.Range("G12").ClearFormats()
